$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: round the Ost/Nord (easting/northing) values to whole numbers.
$ws.Range("Q2").Value = 369470
$ws.Range("R2").Value = 6635346

# Row 3: becomes (a rounded copy of) what used to be row 4's record.
$ws.Range("A3").Value = 111644287
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "35"
$ws.Range("Q3").Value = 369410
$ws.Range("R3").Value = 6635288
$ws.Range("S3").Value = 10
$ws.Range("Z3").Value = "12:37"
$ws.Range("AB3").Value = "12:37"

# Row 4: becomes (a rounded copy of) what used to be row 5's record.
$ws.Range("A4").Value = 111644956
$ws.Range("I4").ClearContents()
$ws.Range("Q4").Value = 369440
$ws.Range("R4").Value = 6635308
$ws.Range("Z4").Value = "12:52"
$ws.Range("AB4").Value = "12:52"

# Row 5: becomes (a rounded copy of) what used to be row 3's record.
$ws.Range("A5").Value = 111644923
$ws.Range("B5").Value = 56543
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 103021
$ws.Range("F5").Value = "Talltita"
$ws.Range("G5").Value = "Poecile montanus"
$ws.Range("H5").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q5").Value = 369436
$ws.Range("R5").Value = 6635294
$ws.Range("S5").Value = 25
